# Apply the "Customer" sales-report edit:
#  - Label cell F8 with "Customer : " and value cell G8 with "Ms. Agoes goes"
#  - Add a new cell comment on G8 ("Customer")
#  - Rename the "Quantity" comments on G11/G12 to "SellQuantity"
#  - Leave the selection on G8 (mirrors the saved selection in the sheet view)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Customer" label/value pair added next to the existing header fields.
$ws.Range("F8").Value = "Customer : "
$ws.Range("G8").Value = "Ms. Agoes goes"

# New comment describing the G8 field.
$null = $ws.Range("G8").AddComment("Customer")

# Existing comments whose wording changed.
$null = $ws.Range("G11").Comment.Text("SellQuantity")
$null = $ws.Range("G12").Comment.Text("SellQuantity")

# Match the saved cursor/selection position.
$null = $ws.Range("G8").Select()
